{"js": "// The edit removes the three bulleted TODO paragraphs:\n//   \"Marcar turnos al cambio\"\n//   \"Verificar email al crearlo\"\n//   \"Verificar email al cambiar de email\"\n// and leaves behind a single empty paragraph (no list/style formatting)\n// that still carries the original \"_GoBack\" bookmark.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Remove the first two list items outright (the third paragraph is kept\n// because it is the one carrying the bookmark that must survive).\nif (paragraphs.items.length > 2) {\n  paragraphs.items[0].delete();\n  paragraphs.items[1].delete();\n  await context.sync();\n}\n\n// Clear the text of the remaining paragraph without disturbing the\n// bookmark that sits after the run (a plain range/paragraph clear() would\n// wipe the bookmark too, so replace the found text instead).\nconst found = body.search(\"Verificar email al cambiar de email\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n\n// Strip the leftover list/style formatting so the paragraph has no pPr,\n// matching a plain default paragraph.\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\n\nif (remaining.items.length > 0) {\n  const last = remaining.items[remaining.items.length - 1];\n  last.detachFromList();\n  last.styleBuiltIn = Word.BuiltInStyleName.normal;\n  await context.sync();\n}\n", "ps1": "# The edit removes the three bulleted TODO paragraphs:\n#   \"Marcar turnos al cambio\"\n#   \"Verificar email al crearlo\"\n#   \"Verificar email al cambiar de email\"\n# and leaves behind a single empty paragraph (no list/style formatting)\n# that still carries the original \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# Delete the first two list items outright. The third paragraph is kept\n# because it is the one that carries the bookmark which must survive.\n$toRemove = @(\"Marcar turnos al cambio\", \"Verificar email al crearlo\")\nforeach ($text in $toRemove) {\n    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Trim() -eq $text) {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n\n# Empty the remaining paragraph's text via Find/Replace (rather than a\n# Range.Delete on the whole paragraph) so the \"_GoBack\" bookmark, which\n# sits right after the run, is left untouched.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n[void]$rng.Find.Execute(\"Verificar email al cambiar de email\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# Strip the leftover list numbering / paragraph style so the paragraph has\n# no pPr left, matching a plain default paragraph.\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$last.Range.ListFormat.RemoveNumbers()\n$last.Style = \"Normal\"\n"}
